$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.832216024398804
$ws.Range("B1").Value = 2.2420814037323
$ws.Range("C1").Value = 2.394904375076294
$ws.Range("D1").Value = 2.866952657699585
$ws.Range("E1").Value = 2.350305080413818
